$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

function PasteFormatsFrom($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---- Row 36 ----
$ws.Range("A36").Value = "PLAUtestuser14"
$ws.Range("B36").Value = "P@ssword2"
$ws.Range("G36").Value = "PLAUtestuser14@mailinator.com"
$ws.Range("J36").Value = "ProdAUtestuser14"
$ws.Range("K36").Value = "ProdAUtestuser14@mailinator.com"

# ---- Row 37 ----
$ws.Range("A37").Value = "PLAUtestuser13"
$ws.Range("B37").Value = "P@ssword2"
$ws.Range("J37").Value = "ProdAUtestuser14"
$ws.Range("K37").Value = "ProdAUtestuser14@mailinator.com"

# ---- Hyperlinks (added before formatting so subsequent format-paste
#      can overwrite the hyperlink's own auto-styling with the exact
#      pre-existing style indices used elsewhere in the sheet) ----
$ws.Hyperlinks.Add($ws.Range("G36"), "mailto:PLAUtestuser14@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B36"), "mailto:P@ssword2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), "mailto:P@ssword2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K36"), "mailto:ProdAUtestuser14@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K37"), "mailto:ProdAUtestuser14@mailinator.com") | Out-Null

# ---- Formatting: copy from comparable existing rows ----
PasteFormatsFrom "A35" "A36"
PasteFormatsFrom "B35" "B36"
PasteFormatsFrom "C35" "C36"
PasteFormatsFrom "D35" "D36"
PasteFormatsFrom "E35" "E36"
PasteFormatsFrom "F35" "F36"
PasteFormatsFrom "G35" "G36"
PasteFormatsFrom "H35" "H36"
PasteFormatsFrom "I35" "I36"
PasteFormatsFrom "J27" "J36"
PasteFormatsFrom "K27" "K36"
PasteFormatsFrom "L35" "L36"

PasteFormatsFrom "A35" "A37"
PasteFormatsFrom "B35" "B37"
PasteFormatsFrom "C35" "C37"
PasteFormatsFrom "D35" "D37"
PasteFormatsFrom "E35" "E37"
PasteFormatsFrom "F35" "F37"
PasteFormatsFrom "C35" "G37"
PasteFormatsFrom "H35" "H37"
PasteFormatsFrom "I35" "I37"
PasteFormatsFrom "J27" "J37"
PasteFormatsFrom "K27" "K37"
PasteFormatsFrom "L35" "L37"

# ---- View state: scroll and selection ----
$ws.Range("D42").Select() | Out-Null
